$d = $word.ActiveDocument

$para = $d.Paragraphs.Last
$r = $para.Range
$r.Text = "Hello WOrld"
